$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 576.41174
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 646.0769
$ws.Range("K17").Value = 1050
$ws.Range("L17").Value = 1938.2307
$ws.Range("M17").Value = -882
$ws.Range("N17").Value = -2274.2307
$ws.Range("H19").Value = 898.7619
$ws.Range("J19").Value = 861.8
$ws.Range("L19").Value = 861.8
$ws.Range("N19").Value = -1211.8
$ws.Range("H32").Value = 2486.2
$ws.Range("J32").Value = 2486.2
$ws.Range("L32").Value = 2486.2
$ws.Range("N32").Value = -3138.2
$ws.Range("H33").Value = 987.6429000000001
$ws.Range("I33").Value = 1173
$ws.Range("K33").Value = 1173
$ws.Range("M33").Value = -944
$ws.Range("H40").Value = 8336128
$ws.Range("I40").Value = 3726.6667
$ws.Range("K40").Value = 3726.6667
$ws.Range("M40").Value = -3551.6667
$ws.Range("H62").Value = 111153500
$ws.Range("I62").Value = 250005250
$ws.Range("K62").Value = 250005250
$ws.Range("M62").Value = -250004626
$ws.Range("H65").Value = 111153500
$ws.Range("I65").Value = 250005250
$ws.Range("K65").Value = 1250026250
$ws.Range("M65").Value = -1250023130
$ws.Range("H100").Value = 2791.9
$ws.Range("I100").Value = 1707.1428
$ws.Range("K100").Value = 1707.1428
$ws.Range("M100").Value = -1166.1428
$ws.Range("H111").Value = 20840666
$ws.Range("I111").Value = 25007404
$ws.Range("J111").Value = 6969
$ws.Range("K111").Value = 75022212
$ws.Range("L111").Value = 20907
$ws.Range("M111").Value = -75019145
$ws.Range("N111").Value = -27041
$ws.Range("H113").Value = 61120652
$ws.Range("I113").Value = 27780754
$ws.Range("J113").Value = 83347250
$ws.Range("K113").Value = 27780754
$ws.Range("L113").Value = 83347250
$ws.Range("M113").Value = -27777500
$ws.Range("N113").Value = -83353758
$ws.Range("H116").Value = 25006748
$ws.Range("I116").Value = 250000000
$ws.Range("K116").Value = 250000000
$ws.Range("M116").Value = -249996558
$ws.Range("H133").Value = 100564.55
$ws.Range("J133").Value = 100564.55
$ws.Range("L133").Value = 100564.55
$ws.Range("N133").Value = -110684.55

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2720570.8
$ws.Range("I32").Value = 2780873
$ws.Range("K32").Value = 2780873
$ws.Range("M32").Value = -2780586
$ws.Range("H132").Value = 4830.109
$ws.Range("I132").Value = 4331.129
$ws.Range("K132").Value = 12993.387
$ws.Range("M132").Value = -10463.387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9261510
$ws.Range("I20").Value = 12822136
$ws.Range("K20").Value = 12822136
$ws.Range("M20").Value = -12821889
$ws.Range("H94").Value = 1532.091
$ws.Range("I94").Value = 500.42856
$ws.Range("K94").Value = 500.42856
$ws.Range("M94").Value = -49.42856
$ws.Range("H138").Value = 83889.75
$ws.Range("J138").Value = 83889.75
$ws.Range("L138").Value = 83889.75
$ws.Range("N138").Value = -94169.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3811.4285
$ws.Range("I16").Value = 2277.9412
$ws.Range("K16").Value = 2277.9412
$ws.Range("M16").Value = -1990.9412
$ws.Range("H31").Value = 6558.9546
$ws.Range("I31").Value = 4811.421
$ws.Range("K31").Value = 4811.421
$ws.Range("M31").Value = -4516.421
$ws.Range("H34").Value = 6558.9546
$ws.Range("I34").Value = 4811.421
$ws.Range("K34").Value = 4811.421
$ws.Range("M34").Value = -4609.421
$ws.Range("H58").Value = 12199863
$ws.Range("I58").Value = 23812002
$ws.Range("K58").Value = 23812002
$ws.Range("M58").Value = -23811799
$ws.Range("H99").Value = 13085.571
$ws.Range("I99").Value = 19266.666
$ws.Range("K99").Value = 19266.666
$ws.Range("M99").Value = -17768.666
$ws.Range("H113").Value = 3811.4285
$ws.Range("I113").Value = 2277.9412
$ws.Range("K113").Value = 2277.9412
$ws.Range("M113").Value = -107.9412000000002
$ws.Range("H120").Value = 73126.336
$ws.Range("J120").Value = 73126.336
$ws.Range("L120").Value = 73126.336
$ws.Range("N120").Value = -80384.336
$ws.Range("H122").Value = 1762
$ws.Range("I122").Value = 1358.75
$ws.Range("K122").Value = 4076.25
$ws.Range("M122").Value = -1626.25
$ws.Range("H126").Value = 13085.571
$ws.Range("I126").Value = 19266.666
$ws.Range("K126").Value = 57799.99800000001
$ws.Range("M126").Value = -55329.99800000001
$ws.Range("H136").Value = 12199863
$ws.Range("I136").Value = 23812002
$ws.Range("K136").Value = 71436006
$ws.Range("M136").Value = -71433456

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3334471
$ws.Range("I5").Value = 4000465.2
$ws.Range("K5").Value = 12001395.6
$ws.Range("M5").Value = -12001283.6
$ws.Range("H68").Value = 3765.7
$ws.Range("I68").Value = 1633.5714
$ws.Range("J68").Value = 8740.666999999999
$ws.Range("K68").Value = 4900.7142
$ws.Range("L68").Value = 26222.001
$ws.Range("M68").Value = -4089.7142
$ws.Range("N68").Value = -27844.001
$ws.Range("H71").Value = 3765.7
$ws.Range("I71").Value = 1633.5714
$ws.Range("J71").Value = 8740.666999999999
$ws.Range("K71").Value = 14702.1426
$ws.Range("L71").Value = 78666.003
$ws.Range("M71").Value = -10646.1426
$ws.Range("N71").Value = -86778.003
$ws.Range("H87").Value = 50011400
$ws.Range("I87").Value = 1000000000
$ws.Range("K87").Value = 3000000000
$ws.Range("M87").Value = -2999998752
$ws.Range("H90").Value = 50011400
$ws.Range("I90").Value = 1000000000
$ws.Range("K90").Value = 9000000000
$ws.Range("M90").Value = -8999993760
$ws.Range("H135").Value = 3334471
$ws.Range("I135").Value = 4000465.2
$ws.Range("K135").Value = 36004186.8
$ws.Range("M135").Value = -36001651.8
$ws.Range("H140").Value = 250643.38
$ws.Range("I140").Value = 250643.38
$ws.Range("K140").Value = 751930.14
$ws.Range("M140").Value = -746750.14

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 16353
$ws.Range("J99").Value = 20000
$ws.Range("L99").Value = 20000
$ws.Range("N99").Value = -24492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 51501.5
$ws.Range("J36").Value = 51501.5
$ws.Range("L36").Value = 51501.5
$ws.Range("N36").Value = -52625.5
$ws.Range("H82").Value = 2029.625
$ws.Range("J82").Value = 2035.875
$ws.Range("L82").Value = 2035.875
$ws.Range("N82").Value = -2757.875
$ws.Range("H85").Value = 2029.625
$ws.Range("J85").Value = 2035.875
$ws.Range("L85").Value = 2035.875
$ws.Range("N85").Value = -4531.875
$ws.Range("H93").Value = 3535.7
$ws.Range("I93").Value = 3200.625
$ws.Range("K93").Value = 3200.625
$ws.Range("M93").Value = -1952.625
$ws.Range("H100").Value = 3222.9565
$ws.Range("I100").Value = 2739.4443
$ws.Range("J100").Value = 3533.7856
$ws.Range("K100").Value = 2739.4443
$ws.Range("L100").Value = 3533.7856
$ws.Range("M100").Value = -2198.4443
$ws.Range("N100").Value = -4615.7856
$ws.Range("H122").Value = 4552.959
$ws.Range("I122").Value = 3899.8235
$ws.Range("J122").Value = 6033.4
$ws.Range("K122").Value = 11699.4705
$ws.Range("L122").Value = 18100.2
$ws.Range("M122").Value = -9249.470499999999
$ws.Range("N122").Value = -23000.2
$ws.Range("H136").Value = 6519.25
$ws.Range("I136").Value = 5401.3687
$ws.Range("J136").Value = 8153.077
$ws.Range("K136").Value = 16204.1061
$ws.Range("L136").Value = 24459.231
$ws.Range("M136").Value = -13654.1061
$ws.Range("N136").Value = -29559.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 55556550
$ws.Range("J107").Value = 55556550
$ws.Range("L107").Value = 166669650
$ws.Range("N107").Value = -166673490
$ws.Range("H132").Value = 20851428
$ws.Range("I132").Value = 23817704
$ws.Range("K132").Value = 71453112
$ws.Range("M132").Value = -71450582
